# 27-09-2023_00 measure unities corrected + notebook created
#
# The "quantidade enzima purificada" column (E) was actually recorded in
# micromolar, not nanomolar, and its real values are 0.2 uM (not the
# placeholder 200 used before). Correct the unit label and the values,
# which ripples into the F column ("degradação mg/diaM" = D/E).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- fix the header labels (E7 / F7) ---------------------------------
$ws.Range("E7").Value = "quantidade enzima purificada (uM)"
$ws.Range("F7").Value = "degradação mg/diaM"

# --- correct the enzyme concentration values (was 200 nM, now 0.2 uM) -
for ($r = 8; $r -le 59; $r++) {
    $ws.Cells.Item($r, 5).Value = 0.2
}

# --- scroll the view back to the top-left corner ----------------------
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("F58").Select()
